$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-measured 3-core linear path traversal timings (column D)
$ws.Range("D3").Value = 0.00183
$ws.Range("D4").Value = 0.54
$ws.Range("D5").Value = 48

# Copy the number formatting/style used by the neighboring column C cells
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("C5").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Move the active selection to D5 as recorded in the saved view state
$ws.Range("D5").Select()
